# Update "Formazione_Franz" sheet: refresh transfer/roster data.
# Rows 3-14 get new player data (transfers in/out, updated roles), and two
# rows (15-16) are appended for players who moved down the list, growing
# the used range from A1:H14 to A1:H16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Davide Zappacosta in (was José Palomino)
$ws.Range("A3").Value = 141312
$ws.Range("B3").Value = "Davide Zappacosta"
$ws.Range("D3").Value = "D, C"

# Row 4: Emil Holm in, moved to Bologna (was Marten De Roon / Atalanta)
$ws.Range("A4").Value = 418579
$ws.Range("B4").Value = "Emil Holm"
$ws.Range("C4").Value = "Bologna"

# Row 5: Rafael Tolói in (was Hans Hateboer)
$ws.Range("A5").Value = 64075
$ws.Range("B5").Value = "Rafael Tolói"
$ws.Range("D5").Value = "D"

# Row 7: Isak Hien in (was Giovanni Bonfanti, moved to row 15)
$ws.Range("A7").Value = 454908
$ws.Range("B7").Value = "Isak Hien"

# Row 8: Luis Muriel moved up (was Ademola Lookman, moved to row 9)
$ws.Range("A8").Value = 84190
$ws.Range("B8").Value = "Luis Muriel"

# Row 9: Ademola Lookman moved down from row 8 (was Teun Koopmeiners)
$ws.Range("A9").Value = 299451
$ws.Range("B9").Value = "Ademola Lookman"
$ws.Range("D9").Value = "C, A"

# Row 10: Marten De Roon moved down from row 4 (was Éderson)
$ws.Range("A10").Value = 85070
$ws.Range("B10").Value = "Marten De Roon"
$ws.Range("D10").Value = "D, C"

# Row 11: Éderson moved down from row 10 (was Leonardo Mendicino)
$ws.Range("A11").Value = 362556
$ws.Range("B11").Value = "Éderson"

# Row 12: Gianluca Scamacca moved up from row 13 (was Luis Muriel)
$ws.Range("A12").Value = 302650
$ws.Range("B12").Value = "Gianluca Scamacca"
$ws.Range("D12").Value = "A"

# Row 13: El Bilal Touré moved up from row 14, now at Stuttgart (was Gianluca Scamacca)
$ws.Range("A13").Value = 391527
$ws.Range("B13").Value = "El Bilal Touré"
$ws.Range("C13").Value = "Stuttgart"

# Row 14: Michel Adopo in, now at Cagliari (was El Bilal Touré)
$ws.Range("A14").Value = 356176
$ws.Range("B14").Value = "Michel Adopo"
$ws.Range("C14").Value = "Cagliari"
$ws.Range("D14").Value = "C"

# Row 15 (new): Giovanni Bonfanti, moved down from row 7
$ws.Range("A15").Value = 432552
$ws.Range("B15").Value = "Giovanni Bonfanti"
$ws.Range("C15").Value = "Atalanta"
$ws.Range("D15").Value = "D"
$ws.Range("E15").Value = "Squadra2"
$ws.Range("F15").Value = "5-4-1"
$ws.Range("G15").Value = ""
$ws.Range("H15").Value = "Franz"

# Row 16 (new): José Palomino, moved down from row 3, now at Cagliari
$ws.Range("A16").Value = 125810
$ws.Range("B16").Value = "José Palomino"
$ws.Range("C16").Value = "Cagliari"
$ws.Range("D16").Value = "D"
$ws.Range("E16").Value = "Squadra2"
$ws.Range("F16").Value = "5-4-1"
$ws.Range("G16").Value = ""
$ws.Range("H16").Value = "Franz"
